$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L5").Value = 3
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 12
$ws.Range("S5").Value = 1.36
$ws.Range("T5").Value = 3
$ws.Range("U5").Value = 1.7
$ws.Range("V5").Value = 2.05
$ws.Range("AC5").Value = 12
$ws.Range("AE5").Value = 13
$ws.Range("AH5").Value = 9
$ws.Range("AT5").Value = 3
$ws.Range("BC5").Value = 501
$ws.Range("P6").Value = 4
$ws.Range("H7").Value = 3.05
$ws.Range("I7").Value = 2.62
$ws.Range("K7").Value = 2.05
$ws.Range("L7").Value = 3.2
$ws.Range("Q7").Value = 1.75
$ws.Range("R7").Value = 1.95
$ws.Range("W7").Value = 10
$ws.Range("X7").Value = 15
$ws.Range("AA7").Value = 20
$ws.Range("AB7").Value = 24
$ws.Range("AC7").Value = 10.5
$ws.Range("AD7").Value = 6.1
$ws.Range("AF7").Value = 45
$ws.Range("AG7").Value = 300
$ws.Range("AH7").Value = 9.75
$ws.Range("AI7").Value = 14.5
$ws.Range("AL7").Value = 21
$ws.Range("AM7").Value = 25
$ws.Range("AN7").Value = 4.65
$ws.Range("AT7").Value = 2.72
$ws.Range("AU7").Value = 6.3
$ws.Range("AV7").Value = 50
$ws.Range("AW7").Value = 4.65
$ws.Range("AY7").Value = 19.5
$ws.Range("AZ7").Value = 60
$ws.Range("BA7").Value = 90
$ws.Range("K8").Value = 2.25
$ws.Range("N8").Value = 12
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 2
$ws.Range("S8").Value = 1.36
$ws.Range("T8").Value = 3
$ws.Range("U8").Value = 1.67
$ws.Range("V8").Value = 2.1
$ws.Range("AE8").Value = 13
$ws.Range("AT8").Value = 3
$ws.Range("AX8").Value = 21
$ws.Range("M10").Value = 1.07
$ws.Range("N10").Value = 9
$ws.Range("G11").Value = 1.95
$ws.Range("I11").Value = 3.9
$ws.Range("L11").Value = 4.75
$ws.Range("AI11").Value = 19
